$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.399.55"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "'2.069.40"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'235.54"
$ws.Range("E5").Value = "  -0.61%  "
$ws.Range("D6").Value = "'0.627"
$ws.Range("E6").Value = "  +1.93%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "'57.49"
$ws.Range("E8").Value = "  -1.03%  "
$ws.Range("E9").Value = "  +3.48%  "
$ws.Range("D10").Value = "'0.0775"
$ws.Range("E10").Value = "  +1.69%  "
$ws.Range("D11").Value = "'0.102"
$ws.Range("E11").Value = "  +0.65%  "
$ws.Range("D12").Value = "'2.375.98"
$ws.Range("E12").Value = "  +0.49%  "
$ws.Range("D13").Value = "'14.46"
$ws.Range("E13").Value = "  -0.14%  "
$ws.Range("D14").Value = "'20.79"
$ws.Range("E14").Value = "  -1.12%  "
$ws.Range("D15").Value = "'0.780"
$ws.Range("E15").Value = "  +0.24%  "
$ws.Range("D16").Value = "'5.18"
$ws.Range("E16").Value = "  -0.31%  "
$ws.Range("D17").Value = "'2.068.32"
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("D18").Value = "'37.351.15"
$ws.Range("E18").Value = "  -0.62%  "
$ws.Range("D19").Value = "'6.31"
$ws.Range("E19").Value = "  +2.13%  "
$ws.Range("D20").Value = "'69.74"
$ws.Range("E20").Value = "  +0.78%  "
$ws.Range("D21").Value = "'0.0₃0818"
$ws.Range("E21").Value = "  +0.39%  "
$ws.Range("D22").Value = "'226.37"
$ws.Range("E22").Value = "  -0.58%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  +1.43%  "
$ws.Range("D25").Value = "'2.41"
$ws.Range("E25").Value = "  -1.48%  "
$ws.Range("D26").Value = "'167.25"
$ws.Range("E26").Value = "  +2.10%  "
$ws.Range("D27").Value = "'8.88"
$ws.Range("E27").Value = "  +0.12%  "
$ws.Range("E28").Value = "  -4.96%  "
$ws.Range("D29").Value = "'19.11"
$ws.Range("E29").Value = "  -0.57%  "
$ws.Range("E30").Value = "  +0.90%  "
$ws.Range("E31").Value = "  -1.31%  "
$ws.Range("D32").Value = "'4.54"
$ws.Range("E32").Value = "  +0.82%  "
$ws.Range("E33").Value = "  -0.62%  "
$ws.Range("E34").Value = "  +1.27%  "
$ws.Range("E35").Value = "  -2.89%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("E37").Value = "  -3.49%  "
$ws.Range("E38").Value = "  +0.08%  "
$ws.Range("E39").Value = "  -4.35%  "
$ws.Range("B40").Value = "Cronos"
$ws.Range("C40").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D40").Value = "'0.0967"
$ws.Range("E40").Value = "  -2.15%  "
$ws.Range("B41").Value = "HuobiToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D41").Value = "'2.94"
$ws.Range("E41").Value = "  -0.66%  "
$ws.Range("D42").Value = "'98.37"
$ws.Range("E42").Value = "  +1.12%  "
$ws.Range("D43").Value = "'1.482.29"
$ws.Range("E43").Value = "  +0.32%  "
$ws.Range("E44").Value = "  +1.04%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").Value = "'4.04"
$ws.Range("E46").Value = "  -12.23%  "
$ws.Range("E47").Value = "  +0.24%  "
$ws.Range("D48").Value = "'15.37"
$ws.Range("E48").Value = "  -3.48%  "
$ws.Range("E49").Value = "  +0.77%  "
$ws.Range("E50").Value = "  +0.91%  "
$ws.Range("D51").Value = "'2.262.72"
$ws.Range("E51").Value = "  +0.59%  "
